# Add a new weekly price record for "Hortaliza, Femacal de La Calera - Choclo".
# This inserts a new row at row 334 (pushing the existing rows 334-364 down to
# 335-365) and populates it with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 334; existing rows 334-364 shift to 335-365.
$ws.Rows.Item(334).Insert()

# Fill in the new row 334 with the new data record.
$ws.Cells.Item(334, 1).Value = 3
$ws.Cells.Item(334, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(334, 3).Value = "Coquimbo"
$ws.Cells.Item(334, 4).Value = 44449
$ws.Cells.Item(334, 5).Value = 5
$ws.Cells.Item(334, 6).Value = 100112024
$ws.Cells.Item(334, 7).Value = "Choclo"
$ws.Cells.Item(334, 8).Value = "Dulce o Americano"
$ws.Cells.Item(334, 9).Value = "Primera"
$ws.Cells.Item(334, 10).Value = 40
$ws.Cells.Item(334, 11).Value = 32000
$ws.Cells.Item(334, 12).Value = 32000
$ws.Cells.Item(334, 13).Value = 32000
$ws.Cells.Item(334, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(334, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(334, 16).Value = 457
$ws.Cells.Item(334, 17).Value = 70
$ws.Cells.Item(334, 18).Value = "Hortaliza"
